# Update balance sheet figures for MLM sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 709000000.0
$ws.Range("C4").Value = 715000000.0
$ws.Range("D4").Value = 713000000.0
$ws.Range("E4").Value = 701000000.0
$ws.Range("F4").Value = 691000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 208000000.0
$ws.Range("C13").Value = 211000000.0
$ws.Range("D13").Value = 197000000.0
$ws.Range("E13").Value = 202000000.0
$ws.Range("F13").Value = 230000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 782000000.0
$ws.Range("C21").Value = 760000000.0
$ws.Range("D21").Value = 741000000.0
$ws.Range("E21").Value = 737000000.0
$ws.Range("F21").Value = 733000000.0
